$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Min Price values
$ws.Range("B2").Value = 3298
$ws.Range("B3").Value = 1399

# Update Date values (numeric serial date, existing cell style already formats as date)
$ws.Range("C2").Value = 44573.549305555556
$ws.Range("C3").Value = 44573.548611111109
